$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 14 so that everything from the old row 14
# onward shifts down by one (rows 14-16 become 15-17), and row 14 is
# free for the new "15:50 / -----" entry inserted between the old
# "Intervalo" row (now 13) and the old "16:40" row (now 15).
$ws.Rows.Item(14).Insert()

# --- Row 8 (11:30): Almoço -> "-" ---
$ws.Range("B8:F8").Value = "-"

# --- Row 9: time 13:00 -> 12:20, values "-" -> Almoço ---
$ws.Range("A9").Value = "12:20"
$ws.Range("B9:F9").Value = "Almoço"

# --- Row 10: time 13:50 -> 13:00 (values stay "-") ---
$ws.Range("A10").Value = "13:00"

# --- Row 11: time 14:40 -> 13:50 (values stay "-") ---
$ws.Range("A11").Value = "13:50"

# --- Row 12: time 15:30 -> 14:40, values Intervalo -> "-" ---
$ws.Range("A12").Value = "14:40"
$ws.Range("B12:F12").Value = "-"

# --- Row 13: time 15:50 -> 15:30, values "-" -> Intervalo ---
$ws.Range("A13").Value = "15:30"
$ws.Range("B13:F13").Value = "Intervalo"

# --- Row 14 (newly inserted): 15:50 / "-" ---
$ws.Range("A14").Value = "15:50"
$ws.Range("B14:F14").Value = "-"

# --- Row 15 (previously row 14): 16:40 / "-" (unchanged values, already shifted) ---
$ws.Range("A15").Value = "16:40"
$ws.Range("B15:F15").Value = "-"

# --- Row 16 (new row): 17:30 / "-" ---
$ws.Range("A16").Value = "17:30"
$ws.Range("B16:F16").Value = "-"

# --- Row 17 (new row): 18:20 / empty ---
$ws.Range("A17").Value = "18:20"
$ws.Range("B17:F17").Value = ""
